$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (values must remain plain text,
# matching the workbook convention where Price/Volume columns store
# numbers and percentages as literal strings, not numeric cells).
$updates = [ordered]@{
    "D2" = "296.87"
    "D3" = "41.89"
    "E3" = "3.63%"
    "D4" = "5.002"
    "E4" = "-0.36%"
    "D5" = "0.07522"
    "E5" = "2.46%"
    "D6" = "4.381"
    "E6" = "1.96%"
    "D7" = "1.580"
    "E7" = "2.89%"
    "D8" = "0.9270"
    "E8" = "0.47%"
    "D10" = "0.1195"
    "E10" = "-1.54%"
    "D11" = "0.1827"
    "E11" = "5.25%"
    "D12" = "0.08925"
    "E12" = "3.53%"
    "E13" = "-5.17%"
    "E14" = "-0.62%"
    "D15" = "0.001286"
    "E15" = "1.50%"
    "D16" = "0.005970"
    "E16" = "0.88%"
    "D19" = "8.084"
    "E19" = "5.29%"
    "D20" = "0.1392"
    "E20" = "0.06%"
    "D22" = "0.04107"
    "E22" = "4.23%"
    "D23" = "0.001267"
    "E23" = "0.42%"
    "D24" = "0.003887"
    "E24" = "2.98%"
    "D25" = "0.0001231"
    "E25" = "-3.91%"
    "D38" = "0.02409"
    "E38" = "5.29%"
    "D39" = "0.05202"
    "E39" = "4.61%"
    "D40" = "0.006308"
    "E40" = "16.72%"
    "D41" = "0.007807"
    "E41" = "1.67%"
    "E42" = "3.22%"
    "D43" = "0.007405"
    "E43" = "0.84%"
    "D44" = "0.007469"
    "E44" = "-6.02%"
    "D45" = "0.2954"
    "E45" = "-6.63%"
    "D46" = "0.00006437"
    "E46" = "1.37%"
    "E48" = "65.77%"
    "D49" = "0.004205"
    "E49" = "0.09%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}

Write-Host "Updated $($updates.Count) cells"
